# Apply "Natmi following Dr Hou advice" edit:
# Expands the Mdk-Itga4 LR-pair sheet from a partial (ECs/FAPs x M2/sCs)
# result into the full 4x4 Sending-cluster x Target-cluster matrix
# (ECs, FAPs, M2, sCs) with updated statistics (16 data rows total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Mdk"
$ws.Cells.Item(2,3).Value = "Itga4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.512729
$ws.Cells.Item(2,8).Value = 7.538187
$ws.Cells.Item(2,9).Value = 0.02190726325199687
$ws.Cells.Item(2,10).Value = 0.02190726325199687
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 23.65990166666667
$ws.Cells.Item(2,14).Value = 70.979705
$ws.Cells.Item(2,15).Value = 0.2997993941754699
$ws.Cells.Item(2,16).Value = 0.29979939417547
$ws.Cells.Item(2,17).Value = 59.45092105498166
$ws.Cells.Item(2,18).Value = 535.058289494835
$ws.Cells.Item(2,19).Value = 0.006567784250991196
$ws.Cells.Item(2,20).Value = 0.006567784250991198

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Mdk"
$ws.Cells.Item(3,3).Value = "Itga4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.512729
$ws.Cells.Item(3,8).Value = 7.538187
$ws.Cells.Item(3,9).Value = 0.02190726325199687
$ws.Cells.Item(3,10).Value = 0.02190726325199687
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.07690566666666666
$ws.Cells.Item(3,14).Value = 0.230717
$ws.Cells.Item(3,15).Value = 0.0009744872400636476
$ws.Cells.Item(3,16).Value = 0.0009744872400636479
$ws.Cells.Item(3,17).Value = 0.1932430988976666
$ws.Cells.Item(3,18).Value = 1.739187890079
$ws.Cells.Item(3,19).Value = 0.0000213483485037862
$ws.Cells.Item(3,20).Value = 0.00002134834850378621

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Mdk"
$ws.Cells.Item(4,3).Value = "Itga4"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.512729
$ws.Cells.Item(4,8).Value = 7.538187
$ws.Cells.Item(4,9).Value = 0.02190726325199687
$ws.Cells.Item(4,10).Value = 0.02190726325199687
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 53.21452433333334
$ws.Cells.Item(4,14).Value = 159.643573
$ws.Cells.Item(4,15).Value = 0.6742919890890982
$ws.Cells.Item(4,16).Value = 0.6742919890890983
$ws.Cells.Item(4,17).Value = 133.7136785135723
$ws.Cells.Item(4,18).Value = 1203.423106622151
$ws.Cells.Item(4,19).Value = 0.01477189211368747
$ws.Cells.Item(4,20).Value = 0.01477189211368748

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Mdk"
$ws.Cells.Item(5,3).Value = "Itga4"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.512729
$ws.Cells.Item(5,8).Value = 7.538187
$ws.Cells.Item(5,9).Value = 0.02190726325199687
$ws.Cells.Item(5,10).Value = 0.02190726325199687
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.967779333333333
$ws.Cells.Item(5,14).Value = 5.903338
$ws.Cells.Item(5,15).Value = 0.02493412949536815
$ws.Cells.Item(5,16).Value = 0.02493412949536816
$ws.Cells.Item(5,17).Value = 4.944496196467333
$ws.Cells.Item(5,18).Value = 44.500465768206
$ws.Cells.Item(5,19).Value = 0.0005462385388144099
$ws.Cells.Item(5,20).Value = 0.0005462385388144101

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Mdk"
$ws.Cells.Item(6,3).Value = "Itga4"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 92.89399466666667
$ws.Cells.Item(6,8).Value = 278.681984
$ws.Cells.Item(6,9).Value = 0.8098976036382196
$ws.Cells.Item(6,10).Value = 0.8098976036382197
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 23.65990166666667
$ws.Cells.Item(6,14).Value = 70.979705
$ws.Cells.Item(6,15).Value = 0.2997993941754699
$ws.Cells.Item(6,16).Value = 0.29979939417547
$ws.Cells.Item(6,17).Value = 2197.862779237191
$ws.Cells.Item(6,18).Value = 19780.76501313472
$ws.Cells.Item(6,19).Value = 0.2428068109149031
$ws.Cells.Item(6,20).Value = 0.2428068109149032

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Mdk"
$ws.Cells.Item(7,3).Value = "Itga4"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 92.89399466666667
$ws.Cells.Item(7,8).Value = 278.681984
$ws.Cells.Item(7,9).Value = 0.8098976036382196
$ws.Cells.Item(7,10).Value = 0.8098976036382197
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.07690566666666666
$ws.Cells.Item(7,14).Value = 0.230717
$ws.Cells.Item(7,15).Value = 0.0009744872400636476
$ws.Cells.Item(7,16).Value = 0.0009744872400636479
$ws.Cells.Item(7,17).Value = 7.144074589169778
$ws.Cells.Item(7,18).Value = 64.296671302528
$ws.Cells.Item(7,19).Value = 0.0007892348805035706
$ws.Cells.Item(7,20).Value = 0.000789234880503571

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Mdk"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 92.89399466666667
$ws.Cells.Item(8,8).Value = 278.681984
$ws.Cells.Item(8,9).Value = 0.8098976036382196
$ws.Cells.Item(8,10).Value = 0.8098976036382197
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 53.21452433333334
$ws.Cells.Item(8,14).Value = 159.643573
$ws.Cells.Item(8,15).Value = 0.6742919890890982
$ws.Cells.Item(8,16).Value = 0.6742919890890983
$ws.Cells.Item(8,17).Value = 4943.309739609871
$ws.Cells.Item(8,18).Value = 44489.78765648884
$ws.Cells.Item(8,19).Value = 0.5461074661157092
$ws.Cells.Item(8,20).Value = 0.5461074661157093

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Mdk"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 92.89399466666667
$ws.Cells.Item(9,8).Value = 278.681984
$ws.Cells.Item(9,9).Value = 0.8098976036382196
$ws.Cells.Item(9,10).Value = 0.8098976036382197
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.967779333333333
$ws.Cells.Item(9,14).Value = 5.903338
$ws.Cells.Item(9,15).Value = 0.02493412949536815
$ws.Cells.Item(9,16).Value = 0.02493412949536816
$ws.Cells.Item(9,17).Value = 182.7948828958436
$ws.Cells.Item(9,18).Value = 1645.153946062592
$ws.Cells.Item(9,19).Value = 0.02019409172710371
$ws.Cells.Item(9,20).Value = 0.02019409172710372

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Mdk"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.610639333333333
$ws.Cells.Item(10,8).Value = 4.831918
$ws.Cells.Item(10,9).Value = 0.0140423817607685
$ws.Cells.Item(10,10).Value = 0.0140423817607685
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 23.65990166666667
$ws.Cells.Item(10,14).Value = 70.979705
$ws.Cells.Item(10,15).Value = 0.2997993941754699
$ws.Cells.Item(10,16).Value = 0.29979939417547
$ws.Cells.Item(10,17).Value = 38.10756824713222
$ws.Cells.Item(10,18).Value = 342.9681142241899
$ws.Cells.Item(10,19).Value = 0.004209897544659065
$ws.Cells.Item(10,20).Value = 0.004209897544659066

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Mdk"
$ws.Cells.Item(11,3).Value = "Itga4"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.610639333333333
$ws.Cells.Item(11,8).Value = 4.831918
$ws.Cells.Item(11,9).Value = 0.0140423817607685
$ws.Cells.Item(11,10).Value = 0.0140423817607685
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.07690566666666666
$ws.Cells.Item(11,14).Value = 0.230717
$ws.Cells.Item(11,15).Value = 0.0009744872400636476
$ws.Cells.Item(11,16).Value = 0.0009744872400636479
$ws.Cells.Item(11,17).Value = 0.1238672916895556
$ws.Cells.Item(11,18).Value = 1.114805625206
$ws.Cells.Item(11,19).Value = 0.0000136841218459714
$ws.Cells.Item(11,20).Value = 0.00001368412184597141

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Mdk"
$ws.Cells.Item(12,3).Value = "Itga4"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.610639333333333
$ws.Cells.Item(12,8).Value = 4.831918
$ws.Cells.Item(12,9).Value = 0.0140423817607685
$ws.Cells.Item(12,10).Value = 0.0140423817607685
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 53.21452433333334
$ws.Cells.Item(12,14).Value = 159.643573
$ws.Cells.Item(12,15).Value = 0.6742919890890982
$ws.Cells.Item(12,16).Value = 0.6742919890890983
$ws.Cells.Item(12,17).Value = 85.70940599589045
$ws.Cells.Item(12,18).Value = 771.384653963014
$ws.Cells.Item(12,19).Value = 0.009468665529017064
$ws.Cells.Item(12,20).Value = 0.009468665529017067

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Mdk"
$ws.Cells.Item(13,3).Value = "Itga4"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.610639333333333
$ws.Cells.Item(13,8).Value = 4.831918
$ws.Cells.Item(13,9).Value = 0.0140423817607685
$ws.Cells.Item(13,10).Value = 0.0140423817607685
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.967779333333333
$ws.Cells.Item(13,14).Value = 5.903338
$ws.Cells.Item(13,15).Value = 0.02493412949536815
$ws.Cells.Item(13,16).Value = 0.02493412949536816
$ws.Cells.Item(13,17).Value = 3.169382793587111
$ws.Cells.Item(13,18).Value = 28.524445142284
$ws.Cells.Item(13,19).Value = 0.0003501345652463976
$ws.Cells.Item(13,20).Value = 0.0003501345652463977

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Mdk"
$ws.Cells.Item(14,3).Value = "Itga4"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 17.68108066666667
$ws.Cells.Item(14,8).Value = 53.04324200000001
$ws.Cells.Item(14,9).Value = 0.154152751349015
$ws.Cells.Item(14,10).Value = 0.154152751349015
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 23.65990166666667
$ws.Cells.Item(14,14).Value = 70.979705
$ws.Cells.Item(14,15).Value = 0.2997993941754699
$ws.Cells.Item(14,16).Value = 0.29979939417547
$ws.Cells.Item(14,17).Value = 418.3326299337345
$ws.Cells.Item(14,18).Value = 3764.99366940361
$ws.Cells.Item(14,19).Value = 0.04621490146491655
$ws.Cells.Item(14,20).Value = 0.04621490146491655

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Mdk"
$ws.Cells.Item(15,3).Value = "Itga4"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 17.68108066666667
$ws.Cells.Item(15,8).Value = 53.04324200000001
$ws.Cells.Item(15,9).Value = 0.154152751349015
$ws.Cells.Item(15,10).Value = 0.154152751349015
$ws.Cells.Item(15,11).Value = 1
$ws.Cells.Item(15,12).Value = 0.3333333333333333
$ws.Cells.Item(15,13).Value = 0.07690566666666666
$ws.Cells.Item(15,14).Value = 0.230717
$ws.Cells.Item(15,15).Value = 0.0009744872400636476
$ws.Cells.Item(15,16).Value = 0.0009744872400636479
$ws.Cells.Item(15,17).Value = 1.359775296057111
$ws.Cells.Item(15,18).Value = 12.237977664514
$ws.Cells.Item(15,19).Value = 0.0001502198892103193
$ws.Cells.Item(15,20).Value = 0.0001502198892103194

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Mdk"
$ws.Cells.Item(16,3).Value = "Itga4"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 17.68108066666667
$ws.Cells.Item(16,8).Value = 53.04324200000001
$ws.Cells.Item(16,9).Value = 0.154152751349015
$ws.Cells.Item(16,10).Value = 0.154152751349015
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 53.21452433333334
$ws.Cells.Item(16,14).Value = 159.643573
$ws.Cells.Item(16,15).Value = 0.6742919890890982
$ws.Cells.Item(16,16).Value = 0.6742919890890983
$ws.Cells.Item(16,17).Value = 940.8902973759631
$ws.Cells.Item(16,18).Value = 8468.012676383667
$ws.Cells.Item(16,19).Value = 0.1039439653306845
$ws.Cells.Item(16,20).Value = 0.1039439653306845

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Mdk"
$ws.Cells.Item(17,3).Value = "Itga4"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 17.68108066666667
$ws.Cells.Item(17,8).Value = 53.04324200000001
$ws.Cells.Item(17,9).Value = 0.154152751349015
$ws.Cells.Item(17,10).Value = 0.154152751349015
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 1.967779333333333
$ws.Cells.Item(17,14).Value = 5.903338
$ws.Cells.Item(17,15).Value = 0.02493412949536815
$ws.Cells.Item(17,16).Value = 0.02493412949536816
$ws.Cells.Item(17,17).Value = 34.79246512686623
$ws.Cells.Item(17,18).Value = 313.132186141796
$ws.Cells.Item(17,19).Value = 0.003843664664203627
$ws.Cells.Item(17,20).Value = 0.003843664664203628
